$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current last row (row 67), shifting the
# "Longest Consecutive Sequence" entry down to row 69, then clone the
# formatting of row 66 onto the two freshly inserted rows so the new
# cells pick up the same styles (s=4,5,6,7,4,6,6,8,4) used throughout
# the table.
$ws.Rows("67:68").Insert()
$ws.Range("A66:I66").Copy()
$ws.Range("A67:I68").PasteSpecial(-4122)

# Row 67 -> problem 66 "Plus One" (加一 / plus-one)
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = "简单"
$ws.Range("C67").Value = "加一"
$ws.Range("D67").Value = "plus-one"
$ws.Range("E67").Value = "200108-1.cpp"
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 8.6
$ws.Range("H67").Value = "100.00%"
$ws.Range("I67").Value = 42668177

# Row 68 -> problem 67 "Add Binary" (二进制求和 / add-binary)
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "简单"
$ws.Range("C68").Value = "二进制求和"
$ws.Range("D68").Value = "add-binary"
$ws.Range("E68").Value = "200108-1.cpp"
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 8.7
$ws.Range("H68").Value = "100.00%"
$ws.Range("I68").Value = 42668665

# The previously-last row (problem 128) now lives at row 69; move the
# active selection there to match the saved workbook view.
$ws.Range("I69").Select()
